$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trans")

# C9: date value 3/6/2019 (serial 43530), formatted like the existing date column
$ws.Range("C9").Value = 43530

# D9 / D10: "PCB" category (reuses existing shared string)
$ws.Range("D9").Value = "PCB"
$ws.Range("D10").Value = "PCB"

# G10: new comment text
$ws.Range("G10").Value = "Cutting line 표시"

# Merge C9:C10 (date now spans both new rows, like the other date groups above)
$ws.Range("C9:C10").Merge()

# Update the active selection to reflect where editing left off
$ws.Range("F10").Select()
